# Update "想去人数" (attendee count) values in the "展览" and "全部类型" sheets.
# Source: gh-pages output regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - row number => new F value
$sheet1Updates = @{
    2  = 38
    3  = 186
    4  = 100
    8  = 1765
    9  = 48
    10 = 13
    11 = 146
    12 = 2016
    13 = 13
    15 = 942
    16 = 455
    17 = 17
    18 = 286
    19 = 207
    23 = 45
    25 = 5
    26 = 1086
    27 = 5
    28 = 325
    29 = 175
    30 = 265
    31 = 313
}

# Sheet "全部类型" (sheet4) - row number => new F value
$sheet4Updates = @{
    2  = 38
    3  = 186
    4  = 100
    8  = 1765
    10 = 48
    11 = 13
    12 = 146
    13 = 2016
    14 = 13
    16 = 943
    17 = 455
    18 = 17
    19 = 286
    20 = 207
    24 = 45
    26 = 5
    27 = 1086
    28 = 5
    29 = 325
    30 = 175
    31 = 265
    32 = 313
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
